$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            }
        }

        $newParts = @()
        if ($hasSystem) {
            $newParts += "System"
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) {
                    $newParts += $p
                }
            }
        } else {
            for ($i = $parts.Count - 1; $i -ge 0; $i--) {
                $newParts += $parts[$i]
            }
        }

        $newVal = $newParts -join ", "
        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
